$wb = $excel.ActiveWorkbook

# The handback report is regenerated: the "Correspond Handoff Datetime" (E)
# and "Correspond Handback DateTime" (H) for the first data row (the
# 26b6cf91... file) move forward to new timestamps in both language sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 16:39:16"
$wsZhCn.Range("H2").Value = "2016-03-12 16:39:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 16:39:19"
$wsDeDe.Range("H2").Value = "2016-03-12 16:39:37"
